$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values per sheet, column F
$updates = @{
    "展览" = @{
        "F2" = 2464
        "F3" = 747
        "F4" = 247
        "F5" = 419
        "F6" = 708
        "F8" = 908
        "F10" = 941
        "F12" = 131
        "F13" = 442
        "F14" = 69
        "F16" = 1099
        "F17" = 24363
        "F18" = 2306
        "F19" = 146
        "F20" = 362
        "F22" = 72
        "F23" = 358
        "F24" = 210
        "F25" = 76
        "F26" = 238
        "F28" = 70
        "F29" = 48
        "F30" = 356
        "F32" = 445
        "F33" = 188
    }
    "演出" = @{
        "F5" = 9
        "F7" = 269
        "F8" = 91
        "F10" = 246
        "F11" = 3636
        "F13" = 154
        "F15" = 17
        "F19" = 136
        "F21" = 4128
    }
    "本地生活" = @{
        "F3" = 170
        "F4" = 790
    }
    "全部类型" = @{
        "F3" = 170
        "F4" = 2464
        "F5" = 790
        "F6" = 747
        "F7" = 247
        "F8" = 419
        "F9" = 708
        "F12" = 9
        "F14" = 269
        "F16" = 908
        "F18" = 941
        "F19" = 131
        "F20" = 442
        "F21" = 69
        "F23" = 1099
        "F24" = 24364
        "F26" = 246
        "F28" = 154
        "F30" = 2306
        "F31" = 146
        "F32" = 17
        "F33" = 362
        "F36" = 358
        "F37" = 210
        "F38" = 238
        "F41" = 70
        "F42" = 48
        "F43" = 136
        "F46" = 445
        "F47" = 188
        "F48" = 4128
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}

$wb.Save()
